$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Apply explicit "Normal" cell style to the existing data range on Sheet1.
# (This creates a dedicated cell-format entry instead of relying on the
# implicit default style index, which is what the fixed column width
# support needs to distinguish formatted ranges from untouched ones.)
$ws1.Range("A2:C4").Style = "Normal"

# Add a second sheet with some numeric sample data + a SUM formula.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = 1
$ws2.Range("A2").Value = 2
$ws2.Range("A3").Value = 3
$ws2.Range("A4").Formula = "=SUM(A1:A3)"

$ws2.Activate()
